# Apply effort-log updates described by the commit:
#  tc05: Timeout condition made much sharper.
#  tc10: Documentation completed
#  doxygen documentation reviewed in parts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Row 47 (2012-11-26): effort reduced from 2 to 1h, comment shortened to
# just "Minor changes on documentation and setup" (drop the appended
# "System load estimation" remark, which now gets its own rows below).
$ws.Range("B47").Value = 1
$ws.Range("D47").Value = "Minor changes on documentation and setup"

# Two new log entries are appended at the end of the table. Column A uses
# the same custom date format/style as the rest of the table, so copy the
# formatting from the previous row (xlPasteFormats = -4122) rather than
# re-typing the number format, which would create a near-duplicate style.
$ws.Range("A49").Copy()
$ws.Range("A50:A51").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A50").Value = 41243
$ws.Range("B50").Value = 3.75
$ws.Range("D50").Value = "tc05 revised, documentation of rtos.c/h extended/corrected"

$ws.Range("A51").Value = 41244
$ws.Range("B51").Value = 1.75
$ws.Range("D51").Value = "Documentation, tc10 and doxygen"

# Restore the view/selection state recorded in the saved workbook.
$ws.Range("E51").Select()
